# Weekly fruit/vegetable price update (Caqui, Macroferia Regional de Talca).
# A new week's record is inserted as row 45, pushing the former rows 45-51
# down to 46-52 (all their data is preserved verbatim by the row insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 45 - shifts old rows 45:51 down to 46:52.
$ws.Rows(45).Insert()

# Populate the new row 45 with this week's record. Columns A,B,C,E,F,G,H,I,J,K
# are identical for every record in this sheet (same market/product), so copy
# them down from the row directly below (now row 46, the old row 45).
$ws.Cells.Item(45, 1).Value = $ws.Cells.Item(46, 1).Value2   # A: Mercado ID
$ws.Cells.Item(45, 2).Value = $ws.Cells.Item(46, 2).Value2   # B: Mercado
$ws.Cells.Item(45, 3).Value = $ws.Cells.Item(46, 3).Value2   # C: Región
$ws.Cells.Item(45, 4).Value = 44722                          # D: Fecha
$ws.Cells.Item(45, 5).Value = $ws.Cells.Item(46, 5).Value2   # E: Codreg
$ws.Cells.Item(45, 6).Value = $ws.Cells.Item(46, 6).Value2   # F: Tipo
$ws.Cells.Item(45, 7).Value = $ws.Cells.Item(46, 7).Value2   # G: Producto ID
$ws.Cells.Item(45, 8).Value = $ws.Cells.Item(46, 8).Value2   # H: Producto
$ws.Cells.Item(45, 9).Value = $ws.Cells.Item(46, 9).Value2   # I: Categoría ID
$ws.Cells.Item(45, 10).Value = $ws.Cells.Item(46, 10).Value2 # J: Categoría
$ws.Cells.Item(45, 11).Value = $ws.Cells.Item(46, 11).Value2 # K: Variedad
$ws.Cells.Item(45, 12).Value = "Primera"                     # L: Calidad
$ws.Cells.Item(45, 13).Value = 25                             # M: Volumen
$ws.Cells.Item(45, 14).Value = 15000                          # N: Precio mínimo
$ws.Cells.Item(45, 15).Value = 15000                          # O: Precio máximo
$ws.Cells.Item(45, 16).Value = 15000                          # P: Precio promedio ponderado
$ws.Cells.Item(45, 17).Value = "$/caja 18 kilos granel"       # Q: Unidad de comercialización
$ws.Cells.Item(45, 18).Value = "Región del Maule"             # R: Origen
$ws.Cells.Item(45, 19).Value = 833                             # S: Precio $/Kg
$ws.Cells.Item(45, 20).Value = 18                              # T: Kg / unidad
